# Apply edit: add survey data through 17 Aug 2020 (rows 193-200)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns AR166:AR168 (value 0) ---
$ws.Range("AR166").Value = 0
$ws.Range("AR167").Value = 0
$ws.Range("AR168").Value = 0

# --- Recomputed values in row 192 ---
$ws.Range("G192").Value = 0.5855221
$ws.Range("AI192").Value = 0.207739
$ws.Range("AW192").Value = 0.8137479
$ws.Range("BB192").Value = 0.4196788

# --- Row 193: fill out the rest of the row (A193 already holds "10 08 2020") ---
# Row 193
$ws.Range("B193").Value = 0.2533333
$ws.Range("C193").Value = 0.9677854
$ws.Range("D193").Value = 0.9594587
$ws.Range("F193").Value = 0.5860857
$ws.Range("G193").Value = 0.5526917
$ws.Range("H193").Value = 0.4741085
$ws.Range("I193").Value = 0.2218402
$ws.Range("J193").Value = 0.2884032
$ws.Range("K193").Value = 0.1827929
$ws.Range("L193").Value = 0.8118982
$ws.Range("M193").Value = 0.8722385
$ws.Range("O193").Value = 0.1934902
$ws.Range("P193").Value = 0.7851693
$ws.Range("Q193").Value = 0.4980506
$ws.Range("R193").Value = 0.4205847
$ws.Range("S193").Value = 0.6903698
$ws.Range("T193").Value = 0.599366
$ws.Range("U193").Value = 0.6695485
$ws.Range("V193").Value = 1.0114849
$ws.Range("W193").Value = 0.269648
$ws.Range("X193").Value = 0.3697466
$ws.Range("Y193").Value = 0.3096986
$ws.Range("Z193").Value = 0.35354
$ws.Range("AA193").Value = 0.4782901
$ws.Range("AB193").Value = 0.5433414
$ws.Range("AD193").Value = 1.011706
$ws.Range("AE193").Value = 0.6755603
$ws.Range("AF193").Value = 0.4868396
$ws.Range("AG193").Value = 0.3121525
$ws.Range("AH193").Value = 0.9308278
$ws.Range("AI193").Value = 0.2841844
$ws.Range("AJ193").Value = 0.2216781
$ws.Range("AK193").Value = 0.5078058
$ws.Range("AL193").Value = 0.824523
$ws.Range("AM193").Value = 0.3328245
$ws.Range("AN193").Value = 0.4718477
$ws.Range("AO193").Value = 0.8155948
$ws.Range("AP193").Value = 0.4733115
$ws.Range("AQ193").Value = 0.4183365
$ws.Range("AS193").Value = 0.3529703
$ws.Range("AT193").Value = 0.6371363
$ws.Range("AU193").Value = 0.4268126
$ws.Range("AV193").Value = 0.8586207
$ws.Range("AW193").Value = 0.8211131
$ws.Range("AX193").Value = 0.7168877
$ws.Range("AY193").Value = 0.3987042
$ws.Range("BA193").Value = 0.3173363
$ws.Range("BB193").Value = 0.415702
$ws.Range("BC193").Value = 0.5220647
$ws.Range("BD193").Value = 0.6944238
$ws.Range("BE193").Value = 0.8915836

# Row 194
$ws.Range("A194").Value = "11 08 2020"
$ws.Range("B194").Value = 0.2764094
$ws.Range("C194").Value = 0.9868111000000001
$ws.Range("D194").Value = 0.9582995
$ws.Range("F194").Value = 0.66605
$ws.Range("G194").Value = 0.5721202
$ws.Range("H194").Value = 0.4710537
$ws.Range("I194").Value = 0.2541164
$ws.Range("J194").Value = 0.2962963
$ws.Range("K194").Value = 0.1856897
$ws.Range("L194").Value = 0.7768072
$ws.Range("M194").Value = 0.7843113
$ws.Range("O194").Value = 0.1873082
$ws.Range("P194").Value = 0.7305501
$ws.Range("Q194").Value = 0.6015201
$ws.Range("R194").Value = 0.4658253
$ws.Range("S194").Value = 0.6596693
$ws.Range("T194").Value = 0.6424748
$ws.Range("U194").Value = 0.6715308
$ws.Range("V194").Value = 0.9944455
$ws.Range("W194").Value = 0.2577387
$ws.Range("X194").Value = 0.3423703
$ws.Range("Y194").Value = 0.274443
$ws.Range("Z194").Value = 0.3707655
$ws.Range("AA194").Value = 0.5212979
$ws.Range("AB194").Value = 0.5408301
$ws.Range("AD194").Value = 0.9161726
$ws.Range("AE194").Value = 0.6111199
$ws.Range("AF194").Value = 0.4603646
$ws.Range("AG194").Value = 0.3466715
$ws.Range("AH194").Value = 0.9933269
$ws.Range("AI194").Value = 0.2501647
$ws.Range("AJ194").Value = 0.2160001
$ws.Range("AK194").Value = 0.5084424
$ws.Range("AL194").Value = 0.8773133
$ws.Range("AM194").Value = 0.3311637
$ws.Range("AN194").Value = 0.4482705
$ws.Range("AO194").Value = 0.7757858
$ws.Range("AP194").Value = 0.4201979
$ws.Range("AQ194").Value = 0.4250822
$ws.Range("AS194").Value = 0.3064956
$ws.Range("AT194").Value = 0.626865
$ws.Range("AU194").Value = 0.433414
$ws.Range("AV194").Value = 0.8664789000000001
$ws.Range("AW194").Value = 0.8294469
$ws.Range("AX194").Value = 0.6366029
$ws.Range("AY194").Value = 0.4429594
$ws.Range("BA194").Value = 0.3064496
$ws.Range("BB194").Value = 0.3818336
$ws.Range("BC194").Value = 0.5039098
$ws.Range("BD194").Value = 0.632173
$ws.Range("BE194").Value = 0.8683086

# Row 195
$ws.Range("A195").Value = "12 08 2020"
$ws.Range("B195").Value = 0.2986757
$ws.Range("C195").Value = 0.994972
$ws.Range("D195").Value = 0.8914474
$ws.Range("F195").Value = 0.641869
$ws.Range("G195").Value = 0.5742868
$ws.Range("H195").Value = 0.5030085
$ws.Range("I195").Value = 0.24567
$ws.Range("J195").Value = 0.3026863
$ws.Range("K195").Value = 0.2060158
$ws.Range("L195").Value = 0.7348446
$ws.Range("M195").Value = 0.868254
$ws.Range("O195").Value = 0.1707007
$ws.Range("P195").Value = 0.661643
$ws.Range("Q195").Value = 0.6552093
$ws.Range("R195").Value = 0.4531449
$ws.Range("S195").Value = 0.6744028
$ws.Range("T195").Value = 0.6674078
$ws.Range("U195").Value = 0.6483692
$ws.Range("V195").Value = 0.9425235
$ws.Range("W195").Value = 0.2707611
$ws.Range("X195").Value = 0.3753943
$ws.Range("Y195").Value = 0.2774169
$ws.Range("Z195").Value = 0.3866196
$ws.Range("AA195").Value = 0.5662248
$ws.Range("AB195").Value = 0.5227462
$ws.Range("AD195").Value = 0.933219
$ws.Range("AE195").Value = 0.5873062
$ws.Range("AF195").Value = 0.4954104
$ws.Range("AG195").Value = 0.4338604
$ws.Range("AH195").Value = 0.8875032
$ws.Range("AI195").Value = 0.2543644
$ws.Range("AJ195").Value = 0.2090971
$ws.Range("AK195").Value = 0.4561835
$ws.Range("AL195").Value = 0.6839326999999999
$ws.Range("AM195").Value = 0.3132796
$ws.Range("AN195").Value = 0.4887202
$ws.Range("AO195").Value = 0.7783167
$ws.Range("AP195").Value = 0.4533425
$ws.Range("AQ195").Value = 0.4534297
$ws.Range("AS195").Value = 0.2657048
$ws.Range("AT195").Value = 0.67862
$ws.Range("AU195").Value = 0.6265219
$ws.Range("AV195").Value = 0.8082175
$ws.Range("AW195").Value = 0.8213824
$ws.Range("AX195").Value = 0.5917479
$ws.Range("AY195").Value = 0.4867049
$ws.Range("BA195").Value = 0.3981358
$ws.Range("BB195").Value = 0.4223112
$ws.Range("BC195").Value = 0.5387677
$ws.Range("BD195").Value = 0.5718478
$ws.Range("BE195").Value = 0.7058314999999999

# Row 196
$ws.Range("A196").Value = "13 08 2020"
$ws.Range("B196").Value = 0.2623991
$ws.Range("C196").Value = 0.998414
$ws.Range("D196").Value = 0.7592411
$ws.Range("F196").Value = 0.6292977
$ws.Range("G196").Value = 0.5801219
$ws.Range("H196").Value = 0.4552175
$ws.Range("I196").Value = 0.2259349
$ws.Range("J196").Value = 0.3051106
$ws.Range("K196").Value = 0.2570456
$ws.Range("L196").Value = 0.7438118
$ws.Range("M196").Value = 0.8965232
$ws.Range("O196").Value = 0.203869
$ws.Range("P196").Value = 0.6634154
$ws.Range("Q196").Value = 0.6481059
$ws.Range("R196").Value = 0.4648753
$ws.Range("S196").Value = 0.6474597
$ws.Range("T196").Value = 0.5775492
$ws.Range("U196").Value = 0.6768379
$ws.Range("V196").Value = 0.9511411
$ws.Range("W196").Value = 0.2530825
$ws.Range("X196").Value = 0.3520323
$ws.Range("Y196").Value = 0.2941067
$ws.Range("Z196").Value = 0.4029897
$ws.Range("AA196").Value = 0.5696909
$ws.Range("AB196").Value = 0.5198987
$ws.Range("AD196").Value = 0.9443015
$ws.Range("AE196").Value = 0.7318568
$ws.Range("AF196").Value = 0.449047
$ws.Range("AG196").Value = 0.4690289
$ws.Range("AH196").Value = 0.878115
$ws.Range("AI196").Value = 0.2547254
$ws.Range("AJ196").Value = 0.2219141
$ws.Range("AK196").Value = 0.4313413
$ws.Range("AL196").Value = 0.7090187999999999
$ws.Range("AM196").Value = 0.3437842
$ws.Range("AN196").Value = 0.5252654
$ws.Range("AO196").Value = 0.7865739
$ws.Range("AP196").Value = 0.4210739
$ws.Range("AQ196").Value = 0.4344144
$ws.Range("AS196").Value = 0.2532907
$ws.Range("AT196").Value = 0.7133047
$ws.Range("AU196").Value = 0.4713272
$ws.Range("AV196").Value = 0.8100362
$ws.Range("AW196").Value = 0.7698126
$ws.Range("AX196").Value = 0.6013991
$ws.Range("AY196").Value = 0.4732275
$ws.Range("BA196").Value = 0.4031291
$ws.Range("BB196").Value = 0.3996197
$ws.Range("BC196").Value = 0.5478204
$ws.Range("BD196").Value = 0.5769803999999999
$ws.Range("BE196").Value = 0.6534792

# Row 197
$ws.Range("A197").Value = "14 08 2020"
$ws.Range("B197").Value = 0.3165138
$ws.Range("C197").Value = 1.0854274
$ws.Range("D197").Value = 0.7002502
$ws.Range("F197").Value = 0.6459673
$ws.Range("G197").Value = 0.5599480999999999
$ws.Range("H197").Value = 0.3754125
$ws.Range("I197").Value = 0.2297438
$ws.Range("J197").Value = 0.251938
$ws.Range("K197").Value = 0.3583174
$ws.Range("L197").Value = 0.705583
$ws.Range("M197").Value = 0.9190444
$ws.Range("O197").Value = 0.1947547
$ws.Range("P197").Value = 0.6795685
$ws.Range("Q197").Value = 0.6717537
$ws.Range("R197").Value = 0.4674359
$ws.Range("S197").Value = 0.6871624
$ws.Range("T197").Value = 0.627196
$ws.Range("U197").Value = 0.6264421999999999
$ws.Range("V197").Value = 0.9422551
$ws.Range("W197").Value = 0.2612901
$ws.Range("X197").Value = 0.3329156
$ws.Range("Y197").Value = 0.2611752
$ws.Range("Z197").Value = 0.3943563
$ws.Range("AA197").Value = 0.5849243
$ws.Range("AB197").Value = 0.5968509
$ws.Range("AD197").Value = 0.9769676
$ws.Range("AE197").Value = 0.5527165000000001
$ws.Range("AF197").Value = 0.4730265
$ws.Range("AG197").Value = 0.4792717
$ws.Range("AH197").Value = 0.8101564
$ws.Range("AI197").Value = 0.2559726
$ws.Range("AJ197").Value = 0.2270597
$ws.Range("AK197").Value = 0.3967432
$ws.Range("AL197").Value = 0.6711781
$ws.Range("AM197").Value = 0.3345428
$ws.Range("AN197").Value = 0.5756863
$ws.Range("AO197").Value = 0.6567484
$ws.Range("AP197").Value = 0.3954896
$ws.Range("AQ197").Value = 0.4335021
$ws.Range("AS197").Value = 0.2567603
$ws.Range("AT197").Value = 0.7663878
$ws.Range("AU197").Value = 0.2731067
$ws.Range("AV197").Value = 0.7819377
$ws.Range("AW197").Value = 0.7426983
$ws.Range("AX197").Value = 0.6499133
$ws.Range("AY197").Value = 0.4695159
$ws.Range("BA197").Value = 0.3892729
$ws.Range("BB197").Value = 0.4259335
$ws.Range("BC197").Value = 0.5100207
$ws.Range("BD197").Value = 0.4978196
$ws.Range("BE197").Value = 0.5598585

# Row 198
$ws.Range("A198").Value = "15 08 2020"
$ws.Range("B198").Value = 0.3295129
$ws.Range("C198").Value = 1.0958915
$ws.Range("D198").Value = 0.6623353
$ws.Range("F198").Value = 0.6450047
$ws.Range("G198").Value = 0.555468
$ws.Range("H198").Value = 0.3899299
$ws.Range("I198").Value = 0.2080339
$ws.Range("J198").Value = 0.2570186
$ws.Range("K198").Value = 0.4053999
$ws.Range("L198").Value = 0.687362
$ws.Range("M198").Value = 0.9430699
$ws.Range("O198").Value = 0.3981009
$ws.Range("P198").Value = 0.6049089
$ws.Range("Q198").Value = 0.6358374
$ws.Range("R198").Value = 0.4766209
$ws.Range("S198").Value = 0.6948769
$ws.Range("T198").Value = 0.5911008
$ws.Range("U198").Value = 0.571309
$ws.Range("V198").Value = 1.0499496
$ws.Range("W198").Value = 0.3266583
$ws.Range("X198").Value = 0.3324555
$ws.Range("Y198").Value = 0.1864575
$ws.Range("Z198").Value = 0.4268539
$ws.Range("AA198").Value = 0.5362803
$ws.Range("AB198").Value = 0.5427921999999999
$ws.Range("AD198").Value = 0.993806
$ws.Range("AE198").Value = 0.6071263
$ws.Range("AF198").Value = 0.4576763
$ws.Range("AG198").Value = 0.5235104
$ws.Range("AH198").Value = 0.7908708
$ws.Range("AI198").Value = 0.2465377
$ws.Range("AJ198").Value = 0.2318847
$ws.Range("AK198").Value = 0.460565
$ws.Range("AL198").Value = 0.8057521
$ws.Range("AM198").Value = 0.3070501
$ws.Range("AN198").Value = 0.5656149
$ws.Range("AO198").Value = 0.6696672
$ws.Range("AP198").Value = 0.3821267
$ws.Range("AQ198").Value = 0.4202235
$ws.Range("AS198").Value = 0.2535806
$ws.Range("AT198").Value = 0.7884158
$ws.Range("AU198").Value = 0.3134851
$ws.Range("AV198").Value = 0.7873492
$ws.Range("AW198").Value = 0.7304414
$ws.Range("AX198").Value = 0.6917541
$ws.Range("AY198").Value = 0.5414987999999999
$ws.Range("BA198").Value = 0.2208477
$ws.Range("BB198").Value = 0.4166878
$ws.Range("BC198").Value = 0.4132926
$ws.Range("BD198").Value = 0.4818452
$ws.Range("BE198").Value = 0.5276143

# Row 199
$ws.Range("A199").Value = "16 08 2020"
$ws.Range("B199").Value = 0.3535354
$ws.Range("C199").Value = 1.0831161
$ws.Range("D199").Value = 0.6764324
$ws.Range("F199").Value = 0.6929037
$ws.Range("G199").Value = 0.5446897000000001
$ws.Range("H199").Value = 0.3929848
$ws.Range("I199").Value = 0.1751877
$ws.Range("J199").Value = 0.2591707
$ws.Range("K199").Value = 0.3399918
$ws.Range("L199").Value = 0.6951433
$ws.Range("M199").Value = 0.9328918
$ws.Range("O199").Value = 0.4293036
$ws.Range("P199").Value = 0.6063935
$ws.Range("Q199").Value = 0.7116460999999999
$ws.Range("R199").Value = 0.4540979
$ws.Range("S199").Value = 0.7939479
$ws.Range("T199").Value = 0.601593
$ws.Range("U199").Value = 0.6104334
$ws.Range("V199").Value = 1.0405785
$ws.Range("W199").Value = 0.2997568
$ws.Range("X199").Value = 0.3764259
$ws.Range("Y199").Value = 0.200538
$ws.Range("Z199").Value = 0.4368032
$ws.Range("AA199").Value = 0.4654381
$ws.Range("AB199").Value = 0.6302356
$ws.Range("AD199").Value = 1.0833467
$ws.Range("AE199").Value = 0.63878
$ws.Range("AF199").Value = 0.4487738
$ws.Range("AG199").Value = 0.6248115
$ws.Range("AH199").Value = 0.6173721
$ws.Range("AI199").Value = 0.2736949
$ws.Range("AJ199").Value = 0.2576313
$ws.Range("AK199").Value = 0.4364019
$ws.Range("AL199").Value = 0.6772361
$ws.Range("AM199").Value = 0.3096298
$ws.Range("AN199").Value = 0.5404335
$ws.Range("AO199").Value = 0.6369735
$ws.Range("AP199").Value = 0.3947716
$ws.Range("AQ199").Value = 0.4498156
$ws.Range("AS199").Value = 0.3840788
$ws.Range("AT199").Value = 0.7512178
$ws.Range("AU199").Value = 0.2991194
$ws.Range("AV199").Value = 0.7413119
$ws.Range("AW199").Value = 0.705906
$ws.Range("AX199").Value = 0.7178339
$ws.Range("AY199").Value = 0.5042001
$ws.Range("BA199").Value = 0.1834859
$ws.Range("BB199").Value = 0.4153139
$ws.Range("BC199").Value = 0.4200087
$ws.Range("BD199").Value = 0.4222172
$ws.Range("BE199").Value = 0.3475469

# --- Row 200: only the date label ---
$ws.Range("A200").Value = "17 08 2020"
